$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 4), mirroring the style/format of the existing rows
$ws.Range("A4").Value = 1905
$ws.Range("C4").Value = 36988
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122) # xlPasteFormats: copy only the number format/style, not the value

# D4 gets the Hyperlink cell style applied but remains empty (no hyperlink added)
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").ClearContents()

# Update selection to match the new active cell
$ws.Range("B4").Select()
